$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Manual Tests")

# Fill in row 8 (Test Number 7) with the new manual test data
$ws.Range("B8").Value = 'I fill every question with "no" or "either", except for StarRating which I say "3" to.'
$ws.Range("C8").Value = 'I am given a lsit of holidays that are 3 starts or higher.'
$ws.Range("D8").Value = 'I am given a lsit of holidays that are 3 starts or higher.'
$ws.Range("E8").Value = "Pass"

# Update the active selection to C13, as in the author's final state
$ws.Range("C13").Select()
